$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend "> " to the pseudo-code lines in column B (B1:B8), matching the
# extra indentation used elsewhere in the sheet to "preserve spacing".
$ws.Range("B1").Value = "> primes = []"
$ws.Range("B2").Value = ">     for a in range(2, 30):"
$ws.Range("B3").Value = ">         for b in range(2, a):"
$ws.Range("B4").Value = ">             if a % b == 0:"
$ws.Range("B5").Value = ">                 break"
$ws.Range("B6").Value = ">     else:"
$ws.Range("B7").Value = ">         primes.append(a)"
$ws.Range("B8").Value = "> print primes"

# Replace the literal booleans in column I with TRUE()/FALSE() formulas and
# drop the custom "TRUE/FALSE" boolean number format back to General.
$ws.Range("I9").Formula = "=FALSE()"
$ws.Range("I9").NumberFormat = "GENERAL"

$ws.Range("I15").Formula = "=TRUE()"
$ws.Range("I15").NumberFormat = "GENERAL"

$ws.Range("I19").Formula = "=FALSE()"
$ws.Range("I19").NumberFormat = "GENERAL"

$ws.Range("I21").Formula = "=FALSE()"
$ws.Range("I21").NumberFormat = "GENERAL"

$ws.Range("I23").Formula = "=FALSE()"
$ws.Range("I23").NumberFormat = "GENERAL"

# Move the active selection from C19 to B10.
$ws.Range("B10").Select()
